# The "median price" row (row 2) is removed from the sheet; the "mean price"
# row (old row 3) and everything below it shifts up by one row, exactly like
# an Excel "Delete Sheet Rows" operation (Rows.Delete shifts cells up and
# carries row-level formatting with the row index, matching the target
# shared-strings / cell-value / dimension changes in the diff).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(2).Delete()
